$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so values like
# "42.235.08", "0.980", "0.0000105" are not coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.235.08"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "2.301.03"
$ws.Range("E3").Value = "  -1.50%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "315.82"
$ws.Range("E5").Value = "  -1.13%  "

# Row 6
$ws.Range("D6").Value = "106.42"
$ws.Range("E6").Value = "  +0.97%  "

# Row 7
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -1.31%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("D10").Value = "40.21"
$ws.Range("E10").Value = "  -0.87%  "

# Row 11
$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12
$ws.Range("D12").Value = "8.41"
$ws.Range("E12").Value = "  +1.21%  "

# Row 13
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14
$ws.Range("D14").Value = "0.980"
$ws.Range("E14").Value = "  -0.44%  "

# Row 15
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16
$ws.Range("D16").Value = "2.651.56"
$ws.Range("E16").Value = "  -1.38%  "

# Row 17
$ws.Range("D17").Value = "2.307.16"
$ws.Range("E17").Value = "  -0.67%  "

# Row 18
$ws.Range("D18").Value = "42.140.89"
$ws.Range("E18").Value = "  -0.88%  "

# Row 19
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").Value = "0.0000105"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$ws.Range("D21").Value = "73.09"
$ws.Range("E21").Value = "  -4.19%  "

# Row 22
$ws.Range("E22").Value = "  -0.50%  "

# Row 23
$ws.Range("D23").Value = "261.01"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").Value = "2.33"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25
$ws.Range("D25").Value = "9.87"
$ws.Range("E25").Value = "  +1.34%  "

# Row 26
$ws.Range("E26").Value = "  +0.54%  "

# Row 27
$ws.Range("D27").Value = "11.02"
$ws.Range("E27").Value = "  -2.63%  "

# Row 28
$ws.Range("E28").Value = "  +2.61%  "

# Row 29
$ws.Range("D29").Value = "22.83"
$ws.Range("E29").Value = "  -0.65%  "

# Row 30
$ws.Range("D30").Value = "36.59"
$ws.Range("E30").Value = "  +2.71%  "

# Row 31
$ws.Range("D31").Value = "166.24"
$ws.Range("E31").Value = "  -4.79%  "

# Row 32
$ws.Range("D32").Value = "0.0895"
$ws.Range("E32").Value = "  +0.89%  "

# Row 33
$ws.Range("E33").Value = "  -1.69%  "

# Row 34
$ws.Range("D34").Value = "5.90"
$ws.Range("E34").Value = "  -2.15%  "

# Row 35
$ws.Range("E35").Value = "  +7.43%  "

# Row 36
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  +2.74%  "

# Row 38
$ws.Range("E38").Value = "  +13.13%  "

# Row 39
$ws.Range("D39").Value = "0.0354"
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -2.78%  "

# Row 41
$ws.Range("D41").Value = "99.58"
$ws.Range("E41").Value = "  +17.89%  "

# Row 42
$ws.Range("E42").Value = "  +1.95%  "

# Row 43
$ws.Range("D43").Value = "71.38"
$ws.Range("E43").Value = "  +2.16%  "

# Row 44
$ws.Range("D44").Value = "0.228"
$ws.Range("E44").Value = "  -0.82%  "

# Row 45
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$ws.Range("D46").Value = "12.35"
$ws.Range("E46").Value = "  +5.03%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "113.62"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "79.67"
$ws.Range("E48").Value = "  +10.12%  "

# Row 49
$ws.Range("D49").Value = "9.19"
$ws.Range("E49").Value = "  +0.87%  "

# Row 50
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").Value = "  -2.50%  "

# Row 51
$ws.Range("E51").Value = "  +3.31%  "
